# Applies the cryptos.xlsx price/volume/ranking update described by the commit:
# "Updated cryptos list on Sun Nov 19 08:24:12 UTC 2023 with GitHub Actions"
#
# Only touches the cells that actually changed (Coin/Link/Price/Volume columns).
# Price-column values that look like plain numbers ("244.81", "0.840", "1.00", ...)
# are forced to Text before the write (NumberFormat "@") so Excel keeps the exact
# source formatting (e.g. trailing zeros) instead of silently turning them into
# numeric values; the format is reset back to Normal right after so no stray style
# survives on cells that do not need one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.708.45"
$ws.Range("E2").Value = "  +0.85%  "

# Row 3
$ws.Range("D3").Value = "1.963.81"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.67%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0799"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.98%  "

# Row 11
$ws.Range("E11").Value = "  +1.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.81%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.87%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.840"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.01%  "

# Row 15
$ws.Range("D15").Value = "2.240.32"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.27%  "

# Row 17
$ws.Range("D17").Value = "1.962.78"
$ws.Range("E17").Value = "  +1.53%  "

# Row 18
$ws.Range("D18").Value = "36.664.03"
$ws.Range("E18").Value = "  +0.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.90%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("E24").Value = "  +6.87%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.144"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.85%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.09%  "

# Row 31
$ws.Range("E31").Value = "  +2.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0619"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.99%  "

# Row 35
$ws.Range("E35").Value = "  +17.64%  "

# Row 36
$ws.Range("E36").Value = "  +6.41%  "

# Row 37
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
$ws.Range("E38").Value = "  -0.38%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.78%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0988"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "

# Row 41
$ws.Range("E41").Value = "  +1.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0212"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.31%  "

# Row 45
$ws.Range("D45").Value = "1.367.25"
$ws.Range("E45").Value = "  +2.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.64%  "

# Row 47
$ws.Range("E47").Value = "  +2.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.52%  "

# Row 49
$ws.Range("E49").Value = "  +0.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.07%  "

# Row 51
$ws.Range("D51").Value = "2.130.28"
$ws.Range("E51").Value = "  +0.98%  "
